# Append: 2025-12-21 18:33 JST
# Update the "取得日時" (retrieved datetime) column (A) for all data rows
# on the active sheet ("ランサーズ") from the previous run timestamp to the
# new one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2:A11").Value = "2025-12-21 18:33:00"
